$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 162.72728
$ws.Range("I5").Value = 91.25
$ws.Range("J5").Value = 353.33334
$ws.Range("K5").Value = 91.25
$ws.Range("L5").Value = 353.33334
$ws.Range("M5").Value = 23.75
$ws.Range("N5").Value = -583.33334

$ws.Range("H47").Value = 19200
$ws.Range("I47").Value = 23000
$ws.Range("J47").Value = 17933.334
$ws.Range("K47").Value = 23000
$ws.Range("L47").Value = 17933.334
$ws.Range("M47").Value = -22028
$ws.Range("N47").Value = -19877.334

$ws.Range("H69").Value = 5004
$ws.Range("I69").Value = 10000
$ws.Range("J69").Value = 4171.3335
$ws.Range("K69").Value = 30000
$ws.Range("L69").Value = 12514.0005
$ws.Range("M69").Value = -29126
$ws.Range("N69").Value = -14262.0005

$ws.Range("H72").Value = 5004
$ws.Range("I72").Value = 10000
$ws.Range("J72").Value = 4171.3335
$ws.Range("K72").Value = 90000
$ws.Range("L72").Value = 37542.0015
$ws.Range("M72").Value = -85632
$ws.Range("N72").Value = -46278.0015

$ws.Range("H132").Value = 1609.079
$ws.Range("I132").Value = 1396.3823
$ws.Range("J132").Value = 3417
$ws.Range("K132").Value = 4189.1469
$ws.Range("L132").Value = 10251
$ws.Range("M132").Value = -1659.1469
$ws.Range("N132").Value = -15311

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 7469.061
$ws.Range("I61").Value = 4093.7693
$ws.Range("K61").Value = 4093.7693
$ws.Range("M61").Value = -3881.7693

$ws.Range("H74").Value = 1815
$ws.Range("I74").Value = 1956.5333
$ws.Range("J74").Value = 1107.3334
$ws.Range("K74").Value = 1956.5333
$ws.Range("L74").Value = 1107.3334
$ws.Range("M74").Value = -1082.5333
$ws.Range("N74").Value = -2855.3334

$ws.Range("H77").Value = 1815
$ws.Range("I77").Value = 1956.5333
$ws.Range("J77").Value = 1107.3334
$ws.Range("K77").Value = 9782.666499999999
$ws.Range("L77").Value = 5536.666999999999
$ws.Range("M77").Value = -5414.666499999999
$ws.Range("N77").Value = -14272.667

$ws.Range("H113").Value = 79058.664
$ws.Range("J113").Value = 79058.664
$ws.Range("L113").Value = 79058.664
$ws.Range("N113").Value = -87736.664

$ws.Range("H132").Value = 1742.7
$ws.Range("I132").Value = 1445.4482
$ws.Range("J132").Value = 2526.3635
$ws.Range("K132").Value = 4336.3446
$ws.Range("L132").Value = 7579.0905
$ws.Range("M132").Value = -1806.3446
$ws.Range("N132").Value = -12639.0905

$ws.Range("H136").Value = 7469.061
$ws.Range("I136").Value = 4093.7693
$ws.Range("K136").Value = 12281.3079
$ws.Range("M136").Value = -9731.3079

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3032169
$ws.Range("I58").Value = 4546663.5
$ws.Range("K58").Value = 4546663.5
$ws.Range("M58").Value = -4546460.5

$ws.Range("H99").Value = 1319.3846
$ws.Range("I99").Value = 1206.3334
$ws.Range("J99").Value = 1416.2858
$ws.Range("K99").Value = 1206.3334
$ws.Range("L99").Value = 1416.2858
$ws.Range("M99").Value = 291.6666
$ws.Range("N99").Value = -4412.2858

$ws.Range("H126").Value = 1319.3846
$ws.Range("I126").Value = 1206.3334
$ws.Range("J126").Value = 1416.2858
$ws.Range("K126").Value = 3619.0002
$ws.Range("L126").Value = 4248.857400000001
$ws.Range("M126").Value = -1149.0002
$ws.Range("N126").Value = -9188.857400000001

$ws.Range("H132").Value = 1991.9354
$ws.Range("J132").Value = 2382.318
$ws.Range("L132").Value = 7146.954000000001
$ws.Range("N132").Value = -12206.954

$ws.Range("H134").Value = 1770.3922
$ws.Range("I134").Value = 1571.5814
$ws.Range("K134").Value = 4714.7442
$ws.Range("M134").Value = -2179.7442

$ws.Range("H136").Value = 3032169
$ws.Range("I136").Value = 4546663.5
$ws.Range("K136").Value = 13639990.5
$ws.Range("M136").Value = -13637440.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()

$ws.Range("H68").Value = 1450
$ws.Range("J68").Value = 1666.6666
$ws.Range("L68").Value = 4999.9998
$ws.Range("N68").Value = -6621.9998

$ws.Range("H71").Value = 1450
$ws.Range("J71").Value = 1666.6666
$ws.Range("L71").Value = 14999.9994
$ws.Range("N71").Value = -23111.9994

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 20031
$ws.Range("J47").Value = 20031
$ws.Range("L47").Value = 20031
$ws.Range("N47").Value = -21167

$ws.Range("H114").Value = 47971.43
$ws.Range("J114").Value = 47971.43
$ws.Range("L114").Value = 47971.43
$ws.Range("N114").Value = -56649.43

$ws.Range("H132").Value = 6683.16
$ws.Range("J132").Value = 10631.692
$ws.Range("L132").Value = 31895.076
$ws.Range("N132").Value = -36955.076

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6622
$ws.Range("I7").Value = 4609.8184
$ws.Range("J7").Value = 14000
$ws.Range("K7").Value = 4609.8184
$ws.Range("L7").Value = 14000
$ws.Range("M7").Value = -4497.8184
$ws.Range("N7").Value = -14224

$ws.Range("H40").Value = 3808.4614
$ws.Range("J40").Value = 5500
$ws.Range("L40").Value = 5500
$ws.Range("N40").Value = -5772

$ws.Range("H126").Value = 6622
$ws.Range("I126").Value = 4609.8184
$ws.Range("J126").Value = 14000
$ws.Range("K126").Value = 13829.4552
$ws.Range("L126").Value = 42000
$ws.Range("M126").Value = -11359.4552
$ws.Range("N126").Value = -46940

$ws.Range("H132").Value = 3750.8215
$ws.Range("I132").Value = 3352.7778
$ws.Range("K132").Value = 10058.3334
$ws.Range("M132").Value = -7528.3334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3856.5
$ws.Range("I62").Value = 3822.4443
$ws.Range("J62").Value = 3890.5557
$ws.Range("K62").Value = 3822.4443
$ws.Range("L62").Value = 3890.5557
$ws.Range("M62").Value = -3198.4443
$ws.Range("N62").Value = -5138.5557

$ws.Range("H65").Value = 3856.5
$ws.Range("I65").Value = 3822.4443
$ws.Range("J65").Value = 3890.5557
$ws.Range("K65").Value = 19112.2215
$ws.Range("L65").Value = 19452.7785
$ws.Range("M65").Value = -15992.2215
$ws.Range("N65").Value = -25692.7785

$ws.Range("H126").Value = 1467.28
$ws.Range("I126").Value = 1510.1052
$ws.Range("J126").Value = 1331.6666
$ws.Range("K126").Value = 4530.3156
$ws.Range("L126").Value = 3994.9998
$ws.Range("M126").Value = -2060.3156
$ws.Range("N126").Value = -8934.9998
